# Rename the dataset id and indicator id on the "metadata" sheet:
#   dataset_internal_id:   LandAndGender      -> LG
#   indicator_internal_id: LandAndGender.3Mb  -> LG.3Mb

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

$ws.Range("B2").Value = "LG"
$ws.Range("B3").Value = "LG.3Mb"
